$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '30.241.20'
$ws.Range('E2').Value = '  -0.87%  '
$ws.Range('D3').Value = '1.877.01'
$ws.Range('E3').Value = '  -1.94%  '
$ws.Range('E4').Value = '  -0.15%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '234.64'
$ws.Range('E5').Value = '  -2.12%  '
$ws.Range('E6').Value = '  -0.12%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4675'
$ws.Range('E7').Value = '  -2.32%  '
$ws.Range('E8').Value = '  -0.60%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06571'
$ws.Range('E9').Value = '  -1.93%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '20.44'
$ws.Range('E10').Value = '  +8.89%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07756'
$ws.Range('E11').Value = '  +0.56%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '97.63'
$ws.Range('E12').Value = '  -4.07%  '
$ws.Range('D13').Value = '1.882.04'
$ws.Range('E13').Value = '  -1.87%  '
$ws.Range('E14').Value = '  -2.59%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.6708'
$ws.Range('E15').Value = '  -0.07%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '283.22'
$ws.Range('E16').Value = '  +7.79%  '
$ws.Range('D17').Value = '30.256.98'
$ws.Range('E17').Value = '  -0.83%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.9995'
$ws.Range('E18').Value = '  -0.23%  '
$ws.Range('B19').Value = 'Avalanche'
$ws.Range('C19').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.59'
$ws.Range('E19').Value = '  -0.61%  '
$ws.Range('B20').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C20').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D20').Value = '2.128.39'
$ws.Range('E20').Value = '  -1.30%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '5.379'
$ws.Range('E21').Value = '  -0.38%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.000007242'
$ws.Range('E22').Value = '  -2.89%  '
$ws.Range('E23').Value = '  -0.18%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '6.166'
$ws.Range('E24').Value = '  -1.98%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '9.333'
$ws.Range('E25').Value = '  -0.18%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '167.67'
$ws.Range('E26').Value = '  +0.48%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '19.14'
$ws.Range('E27').Value = '  +0.10%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.977'
$ws.Range('E28').Value = '  -4.10%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.371'
$ws.Range('E29').Value = '  -1.13%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.09637'
$ws.Range('E30').Value = '  -3.40%  '
$ws.Range('E31').Value = '  -6.45%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.469'
$ws.Range('E32').Value = '  -2.91%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.108'
$ws.Range('E33').Value = '  -3.02%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.04656'
$ws.Range('E34').Value = '  -1.46%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.7006'
$ws.Range('E35').Value = '  -3.35%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.091'
$ws.Range('E36').Value = '  -1.42%  '
$ws.Range('B37').Value = 'HuobiToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.718'
$ws.Range('E37').Value = '  -0.19%  '
$ws.Range('B38').Value = 'VeChain'
$ws.Range('C38').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01859'
$ws.Range('E38').Value = '  -2.85%  '
$ws.Range('B39').Value = 'FraxShare'
$ws.Range('C39').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '6.521'
$ws.Range('E39').Value = '  +4.54%  '
$ws.Range('B40').Value = 'MXToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.523'
$ws.Range('E40').Value = '  -3.57%  '
$ws.Range('B41').Value = 'Aave'
$ws.Range('C41').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '71.94'
$ws.Range('E41').Value = '  -3.63%  '
$ws.Range('B42').Value = 'TrustWalletToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.8617'
$ws.Range('E42').Value = '  -0.11%  '
$ws.Range('B43').Value = 'RenderToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.949'
$ws.Range('E43').Value = '  -0.78%  '
$ws.Range('B44').Value = 'PaxDollar'
$ws.Range('C44').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.001'
$ws.Range('E44').Value = '  -0.12%  '
$ws.Range('B45').Value = 'Quant'
$ws.Range('C45').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '102.98'
$ws.Range('E45').Value = '  -2.33%  '
$ws.Range('B46').Value = 'TheSandbox'
$ws.Range('C46').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.4174'
$ws.Range('E46').Value = '  -1.84%  '
$ws.Range('B47').Value = 'Maker'
$ws.Range('C47').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '981.87'
$ws.Range('E47').Value = '  +5.79%  '
$ws.Range('B48').Value = 'Aptos'
$ws.Range('C48').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '7.206'
$ws.Range('E48').Value = '  -2.69%  '
$ws.Range('B49').Value = 'EnergySwap'
$ws.Range('C49').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '9.143'
$ws.Range('E49').Value = '  +4.44%  '
$ws.Range('B50').Value = 'Elrond'
$ws.Range('C50').Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '33.78'
$ws.Range('E50').Value = '  -2.61%  '
$ws.Range('B51').Value = 'Algorand'
$ws.Range('C51').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.1143'
$ws.Range('E51').Value = '  -4.91%  '
